# Update row 3 of the THREAT_ALERT sheet to reflect the refreshed
# Threat Alert Report data (2026-01-18 01:00 run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("THREAT_ALERT")

$ws.Range("A3").Value = "'22-JAN-26"
$ws.Range("C3").Value = "Nile Air NP-104"
$ws.Range("D3").Value = 388
$ws.Range("E3").Value = 440
$ws.Range("F3").Value = -52
